# Edit script: restructure "Final proposal.docx" per commit
$d = $word.ActiveDocument

# --- Step 1: split the title run into 4 runs and drop "and Gas " ---
$p3 = $d.Paragraphs(3)
$titleRange = $p3.Range
$titleXml = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:b/><w:color w:val=`"FF0000`"/><w:u w:val=`"single`"/></w:rPr><w:t xml:space=`"preserve`">Title: IoT – Based </w:t></w:r><w:r><w:rPr><w:b/><w:color w:val=`"FF0000`"/><w:u w:val=`"single`"/></w:rPr><w:t xml:space=`"preserve`">Air quality </w:t></w:r><w:r><w:rPr><w:b/><w:color w:val=`"FF0000`"/><w:u w:val=`"single`"/></w:rPr><w:t>monitoring</w:t></w:r><w:r><w:rPr><w:b/><w:color w:val=`"FF0000`"/><w:u w:val=`"single`"/></w:rPr><w:t xml:space=`"preserve`"> system.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$d.Range($titleRange.Start, $titleRange.End - 1).InsertXML($titleXml)

# --- Step 2: insert the Objectives/Abstract/Problem-statement block
#             *before* the bookmark paragraph ---
$bookmarkPara = $d.Paragraphs(4)
$beforeRange = $bookmarkPara.Range
$beforeRange.Collapse(1)
$beforeRange.InsertBefore("Objectives`rMain`rTo design a scalable IoT air quality monitor with real-time detection, prediction and alerts.`r`rSpecific`rTo monitor air pollutants such as particulate matter (PM2.5, PM10) and harmful gases (CO, LPG) in real time.`rTo measure environmental parameters like temperature and humidity.`rTo enable remote monitoring through a cloud – connected web application.`rTo provide instant alerts when air quality exceeds safe thresholds.`r`rABSTRACT`r`r`r`rPROBLEM STATEMENT`r`r")

# --- Step 3: insert the Background..References block *after* the
#             bookmark paragraph, anchored on the paragraph that now
#             follows it (re-fetch since indices shifted) ---
$afterAnchor = $d.Paragraphs(21)
$afterRange = $afterAnchor.Range
$afterRange.Collapse(1)
$afterRange.InsertBefore("`rBACKGROUND`r`rLITERATURE REVIEW`r`rMETHODOLOGY`r`r`rEXPECTED OUTCOMES`r`r`r`r`r`r`rRECOMMENDATIONS`r`rREFERENCES`r`r`r`r`r`r")
